# Update "want-to-go" counts (column F) for a handful of events on both
# the "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet,
# which carry duplicated rows for the same events.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 8490
    "F7"  = 1499
    "F18" = 451
    "F20" = 172
    "F23" = 89
    "F24" = 119
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
